$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'88.023.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.11%  "
$ws.Range("D3").Value = "'3.061.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.26%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'209.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.44%  "
$ws.Range("D6").Value = "'617.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.67%  "
$ws.Range("E7").Value = "  -6.08%  "
$ws.Range("D8").Value = "'0.800"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +15.84%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "'3.059.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.29%  "
$ws.Range("D11").Value = "'0.593"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.56%  "
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "'0.0000237"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.94%  "
$ws.Range("D14").Value = "'5.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.30%  "
$ws.Range("D15").Value = "'87.856.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.93%  "
$ws.Range("D16").Value = "'3.630.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.87%  "
$ws.Range("D17").Value = "'31.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.73%  "
$ws.Range("D18").Value = "'3.062.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.84%  "
$ws.Range("D19").Value = "'3.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.98%  "
$ws.Range("D20").Value = "'0.0000200"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -10.91%  "
$ws.Range("D21").Value = "'13.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.95%  "
$ws.Range("D22").Value = "'419.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.88%  "
$ws.Range("D23").Value = "'8.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.96%  "
$ws.Range("D24").Value = "'4.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.49%  "
$ws.Range("D25").Value = "'5.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.69%  "
$ws.Range("D26").Value = "'11.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.66%  "
$ws.Range("D27").Value = "'81.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").Value = "'3.232.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.77%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  +8.73%  "
$ws.Range("D31").Value = "'0.171"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.94%  "
$ws.Range("D32").Value = "'8.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.30%  "
$ws.Range("D33").Value = "'505.70"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.73%  "
$ws.Range("D34").Value = "'3.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -11.38%  "
$ws.Range("D35").Value = "'6.72"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.71%  "
$ws.Range("D36").Value = "'1.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.34%  "
$ws.Range("D37").Value = "'1.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.17%  "
$ws.Range("D38").Value = "'22.21"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.15%  "
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("D40").Value = "'0.129"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "'0.359"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.19%  "
$ws.Range("D44").Value = "'147.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.28%  "
$ws.Range("D45").Value = "'1.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.19%  "
$ws.Range("E46").Value = "  +6.47%  "
$ws.Range("D47").Value = "'43.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.56%  "
$ws.Range("E48").Value = "  +13.84%  "
$ws.Range("D49").Value = "'156.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.83%  "
$ws.Range("E50").Value = "  -5.27%  "
$ws.Range("D51").Value = "'0.699"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.88%  "
